# Updates cryptos list values per the 2024-03-31 GitHub Actions refresh.
# Each target cell is forced to Text format before the write so that
# numeric-looking strings (e.g. "602.24", "70.368.28") are stored as
# text, matching the original inlineStr cells, then the style is reset
# back to Normal so no stray cell formatting is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = '@'
    $rng.Value = $value
    $rng.Style = 'Normal'
}

Set-TextValue 'D2' '70.368.28'
Set-TextValue 'E2' '  +0.78%  '
Set-TextValue 'D3' '3.618.99'
Set-TextValue 'E3' '  +2.38%  '
Set-TextValue 'D4' '0.999'
Set-TextValue 'E4' '  +0.10%  '
Set-TextValue 'D5' '602.24'
Set-TextValue 'E5' '  -0.39%  '
Set-TextValue 'D6' '195.74'
Set-TextValue 'E6' '  -0.33%  '
Set-TextValue 'D7' '0.627'
Set-TextValue 'E7' '  -0.36%  '
Set-TextValue 'D8' '1.00'
Set-TextValue 'E8' '  +0.10%  '
Set-TextValue 'E9' '  +3.76%  '
Set-TextValue 'E10' '  -0.60%  '
Set-TextValue 'D11' '53.23'
Set-TextValue 'E11' '  -0.85%  '
Set-TextValue 'D12' '0.0000305'
Set-TextValue 'E12' '  +0.28%  '
Set-TextValue 'E13' '  +0.64%  '
Set-TextValue 'D14' '4.191.08'
Set-TextValue 'E14' '  +2.45%  '
Set-TextValue 'D15' '598.03'
Set-TextValue 'E15' '  -1.07%  '
Set-TextValue 'D16' '12.94'
Set-TextValue 'E16' '  +1.27%  '
Set-TextValue 'D17' '70.444.27'
Set-TextValue 'E17' '  +0.62%  '
Set-TextValue 'B18' 'Chainlink'
Set-TextValue 'C18' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue 'D18' '19.08'
Set-TextValue 'E18' '  -0.14%  '
Set-TextValue 'B19' 'WrappedEther'
Set-TextValue 'C19' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D19' '3.602.47'
Set-TextValue 'E19' '  +1.51%  '
Set-TextValue 'E20' '  +1.69%  '
Set-TextValue 'D21' '0.999'
Set-TextValue 'E21' '  +0.58%  '
Set-TextValue 'D22' '18.38'
Set-TextValue 'E22' '  +0.64%  '
Set-TextValue 'E23' '  -1.55%  '
Set-TextValue 'D24' '103.14'
Set-TextValue 'E24' '  +0.81%  '
Set-TextValue 'E25' '  -0.23%  '
Set-TextValue 'E26' '  -4.71%  '
Set-TextValue 'D27' '10.63'
Set-TextValue 'E27' '  -2.69%  '
Set-TextValue 'D28' '9.71'
Set-TextValue 'E28' '  +0.62%  '
Set-TextValue 'D29' '33.94'
Set-TextValue 'E29' '  +1.17%  '
Set-TextValue 'D30' '4.70'
Set-TextValue 'E30' '  +8.59%  '
Set-TextValue 'D31' '7.31'
Set-TextValue 'E31' '  +2.66%  '
Set-TextValue 'D32' '12.28'
Set-TextValue 'E32' '  -1.79%  '
Set-TextValue 'E33' '  +1.83%  '
Set-TextValue 'D34' '63.33'
Set-TextValue 'E34' '  +0.23%  '
Set-TextValue 'D35' '0.0₃0891'
Set-TextValue 'E35' '  +1.46%  '
Set-TextValue 'D36' '3.940.72'
Set-TextValue 'E36' '  +5.65%  '
Set-TextValue 'D37' '525.68'
Set-TextValue 'E37' '  +8.12%  '
Set-TextValue 'E38' '  +0.24%  '
Set-TextValue 'D39' '3.07'
Set-TextValue 'E39' '  +0.45%  '
Set-TextValue 'D40' '36.95'
Set-TextValue 'E40' '  +0.90%  '
Set-TextValue 'D41' '0.390'
Set-TextValue 'E41' '  -0.75%  '
Set-TextValue 'D42' '3.54'
Set-TextValue 'E42' '  -2.10%  '
Set-TextValue 'D43' '0.135'
Set-TextValue 'E43' '  +0.98%  '
Set-TextValue 'D44' '0.0462'
Set-TextValue 'E44' '  +1.16%  '
Set-TextValue 'D45' '3.56'
Set-TextValue 'E45' '  +7.77%  '
Set-TextValue 'E46' '  +1.01%  '
Set-TextValue 'E47' '  -0.31%  '
Set-TextValue 'E48' '  -0.08%  '
Set-TextValue 'E49' '  -0.29%  '
Set-TextValue 'D50' '0.000251'
Set-TextValue 'E50' '  -1.84%  '
Set-TextValue 'E51' '  +1.48%  '
